# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
#
# The feed re-associated the odds/result data for several already-listed
# fixtures that share the same match-day block. The row "rank" in column A
# stays put (it's just the sequential row number), but everything from
# column B (match id) through AC (closing Asian-handicap P/L) needs to move
# to a different row within its block. This is a pure cyclic permutation of
# the row payloads, so: snapshot every source row's B:AC values first, then
# write each snapshot into its destination row.
#
# NOTE: this host's PowerShell subset does not bind named parameters
# (e.g. "-Mapping $m" silently yields an empty/default value), so the
# helper below is always invoked positionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Rows {
    param($Mapping)   # hashtable: sourceRow -> destinationRow

    # 1) Snapshot the current B:AC payload for every source row before any
    #    writes happen (the mapping is a cycle, so writing would otherwise
    #    clobber data we still need to read).
    $snapshots = @{}
    foreach ($srcRow in $Mapping.Keys) {
        $snapshots[$srcRow] = $ws.Range("B$srcRow`:AC$srcRow").Value2
    }

    # 2) Push each snapshot into its destination row.
    foreach ($srcRow in $Mapping.Keys) {
        $dstRow = $Mapping[$srcRow]
        $ws.Range("B$dstRow`:AC$dstRow").Value2 = $snapshots[$srcRow]
    }
}

# Block 1 (rows 267-270, match date 2023-11-26): 4-cycle.
Rotate-Rows (@{ 267 = 269; 268 = 270; 269 = 268; 270 = 267 })

# Block 2 (rows 286-291, match date 2023-12-03): 6-cycle.
Rotate-Rows (@{ 286 = 290; 287 = 289; 288 = 291; 289 = 288; 290 = 287; 291 = 286 })
